$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header / account info ---
$ws.Range("C2").Value = "Hartmut"
$ws.Range("B3").Value = "'2570314725427075"
$ws.Range("C3").Value = "Mohaupt"

# --- Opening balance line ---
$ws.Range("D5").Value = "KONTOSTAND AM 13.12.2024"

# --- Transaction row 6 ---
$ws.Range("B6").Value = "17.12."
$ws.Range("C6").Value = "18.12."
$ws.Range("D6").Value = "BURGER KING Duderstadt"
$ws.Range("E6").Value = "41,91-"

# --- Transaction row 7 ---
$ws.Range("B7").Value = "19.12."
$ws.Range("C7").Value = "20.12."
$ws.Range("D7").Value = "ABSCHLAG STROM Stadtwerke Rosenheim 59709452"
$ws.Range("E7").Value = "83,28-"

# --- Transaction row 8 ---
$ws.Range("B8").Value = "23.12."
$ws.Range("C8").Value = "24.12."
$ws.Range("D8").Value = "BEITRAG Allianz SE K-40826349"
$ws.Range("E8").Value = "54,16-"

# --- Transaction row 9: now empty (one fewer transaction than before) ---
$ws.Range("B9:D9").ClearContents()
$ws.Range("E9:F9").ClearContents()
$ws.Range("E9").HorizontalAlignment = -4108
$ws.Range("E9").VerticalAlignment = -4108
$ws.Range("E9").WrapText = $true

# --- Closing balance line ---
$ws.Range("D12").Value = "KONTOSTAND AM 25.12.2024"
$ws.Range("E12").Value = "179,35-"

# --- Next statement date ---
$ws.Range("C13").Value = "IHR NAECHSTER ABRECHNUNGSTERMIN 04.01.2025"
